# Update "想去人数" (number of interested attendees) figures in the
# "展览" and "全部类型" sheets, matching rows by event name so the
# correct cell is updated even though the two sheets have different
# row offsets (全部类型 has one extra row for a concert event).

$wb = $excel.ActiveWorkbook

# Map: event name -> new F-column value
$updates = @{
    "南昌·第四届龙年动漫展——暑假最后的狂欢" = 1031
    "赣州·第五人格only" = 171
    "南昌·Sunflower Garden动漫游戏展" = 2801
    "赣州·卡尼动漫展" = 223
    "鹰潭·MZD动漫游戏嘉年华" = 21
    "南昌·Aud中秋动漫嘉年华" = 123
    "九江·星梦次元XACD动漫游戏博览会国庆盛典" = 72
    "吉安·COMIC LIFE次元假日06" = 84
    "南昌·萌卡动漫展" = 2628
    "江西·JMG（广电）第二届UP动漫游戏博览会" = 807
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count

    for ($r = 2; $r -le $rowCount; $r++) {
        $name = $ws.Cells.Item($r, 3).Value2
        if ($updates.ContainsKey($name)) {
            $ws.Cells.Item($r, 6).Value = $updates[$name]
        }
    }
}
